$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 42
$ws.Range("H42").Value = 99.25
$ws.Range("I42").Value = 86.5
$ws.Range("J42").Value = 137.5
$ws.Range("K42").Value = 259.5
$ws.Range("L42").Value = 412.5
$ws.Range("M42").Value = -29.5
$ws.Range("N42").Value = -872.5

# ALC row 64
$ws.Range("H64").Value = 9500
$ws.Range("J64").Value = 9500
$ws.Range("L64").Value = 9500
$ws.Range("N64").Value = -9996

# ALC row 67
$ws.Range("H67").Value = 9500
$ws.Range("J67").Value = 9500
$ws.Range("L67").Value = 9500
$ws.Range("N67").Value = -11216

# ALC row 70
$ws.Range("H70").Value = 250000000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# ALC row 73
$ws.Range("H73").Value = 250000000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# ALC row 132
$ws.Range("H132").Value = 3472.7932
$ws.Range("I132").Value = 3456.4695
$ws.Range("J132").Value = 3561.6667
$ws.Range("K132").Value = 10369.4085
$ws.Range("L132").Value = 10685.0001
$ws.Range("M132").Value = -7839.408500000001
$ws.Range("N132").Value = -15745.0001

# ALC row 135
$ws.Range("H135").Value = 1884.8
$ws.Range("I135").Value = 1376.6428
$ws.Range("K135").Value = 12389.7852
$ws.Range("M135").Value = -9854.7852

# ALC row 137
$ws.Range("H137").Value = 2280.1
$ws.Range("I137").Value = 1162
$ws.Range("K137").Value = 3486
$ws.Range("M137").Value = -936

# ALC row 138
$ws.Range("H138").Value = 2457.9307
$ws.Range("I138").Value = 1961.3704
$ws.Range("J138").Value = 2755.8667
$ws.Range("K138").Value = 5884.1112
$ws.Range("L138").Value = 8267.6001
$ws.Range("M138").Value = -744.1112000000003
$ws.Range("N138").Value = -18547.6001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 9060216
$ws.Range("I32").Value = 4274807.5
$ws.Range("J32").Value = 35721776
$ws.Range("K32").Value = 4274807.5
$ws.Range("L32").Value = 35721776
$ws.Range("M32").Value = -4274520.5
$ws.Range("N32").Value = -35722350

# ARM row 61
$ws.Range("H61").Value = 2707.28
$ws.Range("I61").Value = 2436.95
$ws.Range("K61").Value = 2436.95
$ws.Range("M61").Value = -2224.95

# ARM row 74
$ws.Range("H74").Value = 2941.45
$ws.Range("I74").Value = 2131.1428
$ws.Range("K74").Value = 2131.1428
$ws.Range("M74").Value = -1257.1428

# ARM row 77
$ws.Range("H77").Value = 2941.45
$ws.Range("I77").Value = 2131.1428
$ws.Range("K77").Value = 10655.714
$ws.Range("M77").Value = -6287.714

# ARM row 122
$ws.Range("H122").Value = 6074.4194
$ws.Range("I122").Value = 4951.244
$ws.Range("K122").Value = 14853.732
$ws.Range("M122").Value = -12403.732

# ARM row 130
$ws.Range("H130").Value = 75003.336
$ws.Range("J130").Value = 75003.336
$ws.Range("L130").Value = 75003.336
$ws.Range("N130").Value = -85043.336

# ARM row 136
$ws.Range("H136").Value = 2707.28
$ws.Range("I136").Value = 2436.95
$ws.Range("K136").Value = 7310.849999999999
$ws.Range("M136").Value = -4760.849999999999

$ws = $wb.Worksheets.Item("BSM")
# BSM row 11
$ws.Range("H11").Value = 2282.4546
$ws.Range("I11").Value = 615.4286
$ws.Range("J11").Value = 5199.75
$ws.Range("K11").Value = 615.4286
$ws.Range("L11").Value = 5199.75
$ws.Range("M11").Value = -475.4286
$ws.Range("N11").Value = -5479.75

# BSM row 94
$ws.Range("H94").Value = 1403.8235
$ws.Range("I94").Value = 867.5833
$ws.Range("J94").Value = 2690.8
$ws.Range("K94").Value = 867.5833
$ws.Range("L94").Value = 2690.8
$ws.Range("M94").Value = -416.5833
$ws.Range("N94").Value = -3592.8

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 3260.423
$ws.Range("I31").Value = 2299.611
$ws.Range("J31").Value = 5422.25
$ws.Range("K31").Value = 2299.611
$ws.Range("L31").Value = 5422.25
$ws.Range("M31").Value = -2004.611
$ws.Range("N31").Value = -6012.25

# CRP row 34
$ws.Range("H34").Value = 3260.423
$ws.Range("I34").Value = 2299.611
$ws.Range("J34").Value = 5422.25
$ws.Range("K34").Value = 2299.611
$ws.Range("L34").Value = 5422.25
$ws.Range("M34").Value = -2097.611
$ws.Range("N34").Value = -5826.25

$ws = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws.Range("H11").Value = 1000370.56
$ws.Range("I11").Value = 1250333.8
$ws.Range("J11").Value = 800400
$ws.Range("K11").Value = 3751001.4
$ws.Range("L11").Value = 2401200
$ws.Range("M11").Value = -3750861.4
$ws.Range("N11").Value = -2401480

# CUL row 26
$ws.Range("H26").Value = 17.75
$ws.Range("I26").Value = 22.333334
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 67.00000199999999
$ws.Range("L26").Value = 12
$ws.Range("M26").Value = 220.999998
$ws.Range("N26").Value = -588

# CUL row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# CUL row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# CUL row 131
$ws.Range("H131").Value = 2091.7693
$ws.Range("I131").Value = 1669.5
$ws.Range("J131").Value = 2168.5454
$ws.Range("K131").Value = 5008.5
$ws.Range("L131").Value = 6505.6362
$ws.Range("M131").Value = 31.5
$ws.Range("N131").Value = -16585.6362

# CUL row 140
$ws.Range("H140").Value = 25001840
$ws.Range("I140").Value = 25001840
$ws.Range("K140").Value = 75005520
$ws.Range("M140").Value = -75000340

$ws = $wb.Worksheets.Item("GSM")
# GSM row 123
$ws.Range("H123").Value = 34668.5
$ws.Range("J123").Value = 34668.5
$ws.Range("L123").Value = 34668.5
$ws.Range("N123").Value = -39568.5

# GSM row 134
$ws.Range("H134").Value = 86663
$ws.Range("J134").Value = 86663
$ws.Range("L134").Value = 259989
$ws.Range("N134").Value = -265059

$ws = $wb.Worksheets.Item("LTW")
# LTW row 43
$ws.Range("H43").Value = 142866430
$ws.Range("I43").Value = 7991
$ws.Range("J43").Value = 166676160
$ws.Range("K43").Value = 7991
$ws.Range("L43").Value = 166676160
$ws.Range("M43").Value = -7798
$ws.Range("N43").Value = -166676546

# LTW row 55
$ws.Range("H55").Value = 305.32257
$ws.Range("I55").Value = 200.72223
$ws.Range("K55").Value = 200.72223
$ws.Range("M55").Value = -27.72223

# LTW row 82
$ws.Range("H82").Value = 3763.611
$ws.Range("I82").Value = 3094.2222
$ws.Range("K82").Value = 3094.2222
$ws.Range("M82").Value = -2733.2222

# LTW row 85
$ws.Range("H85").Value = 3763.611
$ws.Range("I85").Value = 3094.2222
$ws.Range("K85").Value = 3094.2222
$ws.Range("M85").Value = -1846.2222

# LTW row 94
$ws.Range("H94").Value = 49999
$ws.Range("J94").Value = 49999
$ws.Range("L94").Value = 49999
$ws.Range("N94").Value = -51351

# LTW row 135
$ws.Range("H135").Value = 49993.5
$ws.Range("J135").Value = 49993.5
$ws.Range("L135").Value = 49993.5
$ws.Range("N135").Value = -60133.5

# LTW row 136
$ws.Range("H136").Value = 3540.4849
$ws.Range("I136").Value = 3140.5
$ws.Range("K136").Value = 9421.5
$ws.Range("M136").Value = -6871.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 4530.304
$ws.Range("I122").Value = 3794.5881
$ws.Range("J122").Value = 6614.8335
$ws.Range("K122").Value = 11383.7643
$ws.Range("L122").Value = 19844.5005
$ws.Range("M122").Value = -8933.764299999999
$ws.Range("N122").Value = -24744.5005

# WVR row 132
$ws.Range("H132").Value = 3292.4375
$ws.Range("I132").Value = 2745.3462
$ws.Range("K132").Value = 8236.0386
$ws.Range("M132").Value = -5706.0386
